$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("README")
try {
    $st = $wb.Styles.Add("MyTitle")
    Write-Host "added style"
    $st.Font.Bold = $true
    $st.Font.Size = 16
    $st.Font.Name = "Calibri Light"
    $ws.Range("A6").Style = "MyTitle"
    Write-Host "applied"
} catch {
    Write-Host ("ERR: " + $_)
}
